# Add support for custom run names:
# Insert a new row in the settings table (row 13) with a new
# "Run Names" field (RUN.NAME_FORMAT), default value
# "[FuncName]-S1R[R#1]" and an explanatory comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 13 (a blank spacer row).
# This shifts every following row down by one, which matches the
# rest of the sheet (table, dimension, etc.) expanding by one row.
[void]$ws.Rows.Item(13).Insert()

# Fill in the new row's contents.
$ws.Range("B13").Value = "Run Names"
$ws.Range("C13").Value = "[FuncName]-S1R[R#1]"

$commentPart1 = 'When entered correctly, BV should name runs in the format of "[FuncName]-S1R[R#1]". '
$commentPart2 = "This field accomodates other run naming schemes, but can typically be left as the default."
$ws.Range("D13").Value = $commentPart1 + $commentPart2
$len1 = $commentPart1.Length
$totalLen = ($commentPart1 + $commentPart2).Length
$ws.Range("D13").Characters($len1 + 1, $totalLen - $len1).Font.Bold = $true

$ws.Range("E13").Value = "RUN.NAME_FORMAT"

# Match row height used by the other multi-line comment rows.
$ws.Rows.Item(13).RowHeight = 45

# Grow the table (ListObject) to include the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lastRow = $lo.Range.Rows.Count + 1
[void]$lo.Resize($ws.Range("A1:E" + $lastRow))

# Column E now needs to fit the new "RUN.NAME_FORMAT" value.
[void]$ws.Columns.Item(5).AutoFit()

# Scroll back to the default top-left position and select the new cell.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E13").Select()
